$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("measures")

# Insert a new row at position 7, shifting existing rows 7+ down to 8+.
$ws.Rows(7).Insert()

# Populate the new row 7 with the "Winner concentration" measure.
$ws.Range("A7").Value = "Winner concentration"
$ws.Range("B7").Value = "Herfindahl-Hirschman index of winners' share of agency budget"
$ws.Range("E7").Value = "Anti-trust literature"
$ws.Range("K7").Value = "Draft"

# Match formatting of surrounding measure rows: column B uses wrap-text style 3.
$ws.Range("B7").WrapText = $true

# Update the "Repeat winners" rows (5 and 6) status from Draft to Draft - not used.
$ws.Range("K5").Value = "Draft - not used"
$ws.Range("K6").Value = "Draft - not used"

# Restore view state: pane frozen at row1/col1, scrolled to top, with F7 selected.
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("F7").Select()
